$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.415.64'
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("D3").Value = '1.960.08'
$ws.Range("E3").Value = '  -4.52%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.77%  '
$ws.Range("E6").Value = '  -4.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.71'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -9.08%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.372'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '55.63'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.23%  '
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.19'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.830'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -8.10%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.247.06'
$ws.Range("E15").Value = '  -4.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.48'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -7.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.31'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.98%  '
$ws.Range("D18").Value = '1.967.71'
$ws.Range("E18").Value = '  -4.19%  '
$ws.Range("D19").Value = '36.335.84'
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.52'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.10%  '
$ws.Range("D21").Value = '0.0₃0879'
$ws.Range("E21").Value = '  -1.43%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.13'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.84%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.14'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.40%  '
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '166.06'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.73'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("E30").Value = '  -10.95%  '
$ws.Range("E31").Value = '  -3.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.16'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.74'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0636'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.32'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.43%  '
$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.13'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.48%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.16'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -9.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0969'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.18'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -6.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0210'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.74'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.05'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.81'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.85%  '
$ws.Range("D48").Value = '1.339.39'
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.24'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -7.12%  '
$ws.Range("E50").Value = '  -3.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.77'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.53%  '
